$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = 0.8764979839324951
$ws.Range("C2").Value = 0.5605244040489197
$ws.Range("D2").Value = 0.8648582696914673
$ws.Range("E2").Value = 0.8836954236030579
$ws.Range("F2").Value = 0.8707399964332581
$ws.Range("G2").Value = 86.64286804199219
$ws.Range("H2").Value = 8.461791038513184
$ws.Range("I2").Value = 8.023092269897461
$ws.Range("J2").Value = 86.33681488037109
$ws.Range("K2").Value = 88.83335876464844
